# Applies the weekly Fruta/hortaliza reshuffle for "Mora" rows (Vega Central
# Mapocho de Santiago): dates, volumes, prices and origin get realigned to
# a different row ordering, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44533
$ws.Cells.Item(2, 13).Value = 150
$ws.Cells.Item(2, 14).Value = 4000
$ws.Cells.Item(2, 15).Value = 4000
$ws.Cells.Item(2, 16).Value = 4000
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 2000

# Row 3
$ws.Cells.Item(3, 4).Value = 44978
$ws.Cells.Item(3, 13).Value = 500
$ws.Cells.Item(3, 14).Value = 3000
$ws.Cells.Item(3, 15).Value = 3000
$ws.Cells.Item(3, 16).Value = 3000
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 1500

# Row 4
$ws.Cells.Item(4, 4).Value = 44194
$ws.Cells.Item(4, 13).Value = 120
$ws.Cells.Item(4, 14).Value = 3000
$ws.Cells.Item(4, 15).Value = 3000
$ws.Cells.Item(4, 16).Value = 3000
$ws.Cells.Item(4, 18).Value = "Provincia de Linares"
$ws.Cells.Item(4, 19).Value = 1500

# Row 5
$ws.Cells.Item(5, 4).Value = 44174
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 3200
$ws.Cells.Item(5, 15).Value = 3200
$ws.Cells.Item(5, 16).Value = 3200
$ws.Cells.Item(5, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(5, 19).Value = 1600

# Row 6
$ws.Cells.Item(6, 4).Value = 44208
$ws.Cells.Item(6, 13).Value = 85
$ws.Cells.Item(6, 14).Value = 3000
$ws.Cells.Item(6, 15).Value = 3000
$ws.Cells.Item(6, 16).Value = 3000
$ws.Cells.Item(6, 18).Value = "Provincia de Linares"
$ws.Cells.Item(6, 19).Value = 1500

# Row 8
$ws.Cells.Item(8, 4).Value = 44582
$ws.Cells.Item(8, 13).Value = 380
$ws.Cells.Item(8, 14).Value = 5000
$ws.Cells.Item(8, 15).Value = 5000
$ws.Cells.Item(8, 16).Value = 5000
$ws.Cells.Item(8, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(8, 19).Value = 2500

# Row 9
$ws.Cells.Item(9, 4).Value = 44236
$ws.Cells.Item(9, 13).Value = 300
$ws.Cells.Item(9, 14).Value = 3600
$ws.Cells.Item(9, 15).Value = 4000
$ws.Cells.Item(9, 16).Value = 3800
$ws.Cells.Item(9, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(9, 19).Value = 1900

# Row 10
$ws.Cells.Item(10, 4).Value = 44980
$ws.Cells.Item(10, 13).Value = 250
$ws.Cells.Item(10, 14).Value = 4000
$ws.Cells.Item(10, 15).Value = 4000
$ws.Cells.Item(10, 16).Value = 4000
$ws.Cells.Item(10, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(10, 19).Value = 2000

# Row 11
$ws.Cells.Item(11, 4).Value = 44232
$ws.Cells.Item(11, 13).Value = 200
$ws.Cells.Item(11, 14).Value = 3000
$ws.Cells.Item(11, 15).Value = 3000
$ws.Cells.Item(11, 16).Value = 3000
$ws.Cells.Item(11, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(11, 19).Value = 1500

# Row 12
$ws.Cells.Item(12, 4).Value = 44188
$ws.Cells.Item(12, 13).Value = 150
$ws.Cells.Item(12, 14).Value = 3000
$ws.Cells.Item(12, 15).Value = 3400
$ws.Cells.Item(12, 16).Value = 3240
$ws.Cells.Item(12, 18).Value = "Provincia de Linares"
$ws.Cells.Item(12, 19).Value = 1620

# Row 13
$ws.Cells.Item(13, 4).Value = 44238
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 3600
$ws.Cells.Item(13, 15).Value = 4000
$ws.Cells.Item(13, 16).Value = 3800
$ws.Cells.Item(13, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(13, 19).Value = 1900

# Row 14
$ws.Cells.Item(14, 4).Value = 44168
$ws.Cells.Item(14, 13).Value = 170
$ws.Cells.Item(14, 14).Value = 8000
$ws.Cells.Item(14, 15).Value = 8000
$ws.Cells.Item(14, 16).Value = 8000
$ws.Cells.Item(14, 18).Value = "Provincia de Linares"
$ws.Cells.Item(14, 19).Value = 4000

# Row 15
$ws.Cells.Item(15, 4).Value = 44617
$ws.Cells.Item(15, 13).Value = 90
$ws.Cells.Item(15, 14).Value = 6500
$ws.Cells.Item(15, 15).Value = 6500
$ws.Cells.Item(15, 16).Value = 6500
$ws.Cells.Item(15, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(15, 19).Value = 3250

# Row 16
$ws.Cells.Item(16, 4).Value = 44231
$ws.Cells.Item(16, 13).Value = 150
$ws.Cells.Item(16, 14).Value = 3400
$ws.Cells.Item(16, 15).Value = 3400
$ws.Cells.Item(16, 16).Value = 3400
$ws.Cells.Item(16, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(16, 19).Value = 1700

# Row 17
$ws.Cells.Item(17, 4).Value = 44237
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 3600
$ws.Cells.Item(17, 15).Value = 4000
$ws.Cells.Item(17, 16).Value = 3800
$ws.Cells.Item(17, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(17, 19).Value = 1900
